# Fruta / hortaliza, semanal
# Insert one new weekly record at row 234 of the "Perejil" (Vega Modelo de
# Temuco) data sheet. Excel's native row-insert shifts rows 234-290 down to
# 235-291 (and the rest of the shared columns tag along automatically), so
# we only need to populate the newly-opened row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 234..290 down to 235..291, opening up a blank row 234.
$ws.Rows(234).Insert()

# Populate the new row 234 with the new weekly price record.
$ws.Cells.Item(234, 1).Value = 10
$ws.Cells.Item(234, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(234, 3).Value = "La Araucanía"
$ws.Cells.Item(234, 4).Value = 44642
$ws.Cells.Item(234, 5).Value = 9
$ws.Cells.Item(234, 6).Value = 100112044
$ws.Cells.Item(234, 7).Value = "Perejil"
$ws.Cells.Item(234, 8).Value = "Sin especificar"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 35
$ws.Cells.Item(234, 11).Value = 5000
$ws.Cells.Item(234, 12).Value = 5000
$ws.Cells.Item(234, 13).Value = 5000
$ws.Cells.Item(234, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(234, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(234, 16).Value = 1667
$ws.Cells.Item(234, 17).Value = 3
$ws.Cells.Item(234, 18).Value = "Hortaliza"
